$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# 1. Update the header metadata strings (shared between both sheets).
#    Row 2 = report title, Row 6 = release date, Row 7 = next update.
# -------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Item(2, 1).Value = "June 2018 Monthly Energy Review"
    $ws.Cells.Item(6, 1).Value = "Release Date: June 26, 2018"
    $ws.Cells.Item(7, 1).Value = "Next Update: July 26, 2018"
}

# -------------------------------------------------------------------
# 2. "Monthly Data" sheet: revise Jan-2016..Feb-2017 (rows 529-542)
#    with refreshed EIA figures, then append Mar-2017..Mar-2018
#    (rows 543-555) as brand-new rows.
# -------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Data")

$monthlyUpdates = @(
    @(529, 75.262, 74.912, 150.174, 1328.483, 1397.261, 1651.675, 3048.936, 4377.419, 0, 62134.631, 66662.224),
    @(530, 75.387, 75.036, 150.423, 1361.445, 1282.036, 1755.363, 3037.399, 4398.844, 0, 50661.45,  55210.717),
    @(531, 74.003, 73.658, 147.661, 1433.657, 1274.833, 1770.31,  3045.143, 4478.8,   0, 39948.145, 44574.606),
    @(532, 45.553, 28.657, 74.21,   1324.031, 1075.583, 1750.917, 2826.5,   4150.531, 0, 39158.963, 43383.704),
    @(533, 36.542, 22.989, 59.531,  1366.87,  1178.029, 1656.568, 2834.597, 4201.467, 0, 45081.935, 49342.933),
    @(534, 46.166, 29.043, 75.209,  1404.818, 1242.702, 1578.086, 2820.788, 4225.606, 0, 63250.414, 67551.229),
    @(535, 46.143, 17.383, 63.526,  1432.54,  1321.161, 1514.584, 2835.745, 4268.285, 0, 74236.728, 78568.539),
    @(536, 49.413, 18.615, 68.028,  1394.678, 1292.162, 1529.738, 2821.9,   4216.578, 0, 73889.93,  78174.536),
    @(537, 49.607, 18.688, 68.295,  1335.799, 1157.257, 1668.331, 2825.588, 4161.387, 0, 62385.216, 66614.898),
    @(538, 49.856, 37.991, 87.847,  1334.67,  1126.246, 1782.495, 2908.741, 4243.411, 0, 54621.445, 58952.703),
    @(539, 59.538, 45.368, 104.906, 1325.968, 1092.706, 1830.459, 2923.165, 4249.133, 0, 48179.203, 52533.242),
    @(540, 75.425, 57.474, 132.899, 1441.749, 1280.244, 1640.041, 2920.285, 4362.034, 0, 65006.425, 69501.358),
    @(541, 66.184, 72.161, 138.345, 1430.645, 1290.391, 1553.562, 2843.953, 4274.598, 0, 63547.714, 67960.657),
    @(542, 53.674, 58.521, 112.195, 1367.727, 1087.427, 1766.993, 2854.42,  4222.147, 0, 47964.848, 52299.19)
)

foreach ($row in $monthlyUpdates) {
    $r = $row[0]
    for ($col = 3; $col -le 13; $col++) {
        $wsMonthly.Cells.Item($r, $col).Value = $row[$col - 2]
    }
}

$monthlyNewRows = @(
    @(543, 42795, 58.423, 63.699, 122.122, 1437.669, 1172.172, 1664.293, 2836.465, 4274.134, 0, 48825.958, 53222.214),
    @(544, 42826, 40.473, 24.9,   65.373,  1440.81,  1067.605, 1629.855, 2697.46,  4138.27,  0, 44323.847, 48527.49),
    @(545, 42856, 39.962, 24.586, 64.548,  1482.486, 1098.283, 1604.724, 2703.007, 4185.493, 0, 50926.005, 55176.046),
    @(546, 42887, 45.507, 27.997, 73.504,  1401.664, 1094.108, 1617.146, 2711.254, 4112.918, 0, 58951.924, 63138.346),
    @(547, 42917, 53.309, 16.546, 69.855,  1494.46,  1047.123, 1838.301, 2885.424, 4379.884, 0, 69900.111, 74349.85),
    @(548, 42948, 48.549, 15.069, 63.618,  1528.056, 1064.994, 1807.254, 2872.248, 4400.304, 0, 65933.994, 70397.916),
    @(549, 42979, 47.069, 14.609, 61.678,  1468.767, 1030.015, 1809.249, 2839.264, 4308.031, 0, 54779.784, 59149.493),
    @(550, 43009, 42.669, 37.703, 80.372,  1469.57,  1149.103, 1641.732, 2790.835, 4260.405, 0, 50214.467, 54555.244),
    @(551, 43040, 49.578, 43.808, 93.386,  1456.863, 1142.217, 1650.271, 2792.488, 4249.351, 0, 50992.13,  55334.867),
    @(552, 43070, 61.661, 54.486, 116.147, 1558.946, 1180.736, 1605.369, 2786.105, 4345.051, 0, 58388.345, 62849.543),
    @(553, 43101, 68.861, 35.059, 103.92,  1689.106, 1269.77,  1667.075, 2936.845, 4625.951, 0, 64650.176, 69380.047),
    @(554, 43132, 53.112, 50.81,  103.922, 1388.187, 1131.822, 1729.148, 2860.97,  4249.157, 0, 45823.067, 50176.146),
    @(555, 43160, 50.608, 7.14,   57.748,  1113.576, 1169.107, 1610.478, 2779.585, 3893.161, 0, 44495.503, 48446.411)
)

foreach ($row in $monthlyNewRows) {
    $r = $row[0]

    $wsMonthly.Cells.Item($r, 1).NumberFormat = "yyyy mmmm"
    $wsMonthly.Cells.Item($r, 1).HorizontalAlignment = -4131
    $wsMonthly.Cells.Item($r, 1).Value = $row[1]

    $wsMonthly.Cells.Item($r, 2).Value = "Not Available"

    for ($col = 3; $col -le 13; $col++) {
        $wsMonthly.Cells.Item($r, $col).Value = $row[$col - 1]
    }
}

# -------------------------------------------------------------------
# 3. "Annual Data" sheet: revise the 2016 total (row 80) and append
#    the new 2017 annual total as row 81.
# -------------------------------------------------------------------
$wsAnnual = $wb.Worksheets.Item("Annual Data")

$annualRow80 = @(682.895, 499.814, 1182.709, 16484.708, 14720.22, 20128.567, 34848.787, 51333.495, 0, 678554.486, 731070.69)
for ($col = 3; $col -le 13; $col++) {
    $wsAnnual.Cells.Item(80, $col).Value = $annualRow80[$col - 3]
}

$wsAnnual.Cells.Item(81, 1).HorizontalAlignment = -4131
$wsAnnual.Cells.Item(81, 1).Value = 2017
$wsAnnual.Cells.Item(81, 2).Value = "Not Available"

$annualRow81 = @(607.058, 454.085, 1061.143, 17537.663, 13424.174, 20188.749, 33612.923, 51150.586, 0, 664749.129, 716960.858)
for ($col = 3; $col -le 13; $col++) {
    $wsAnnual.Cells.Item(81, $col).Value = $annualRow81[$col - 3]
}
